$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 65

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null -and $val.Length -gt 0) {
        $first = $val.Substring(0, 1).ToLower()
        $rest = $val.Substring(1)
        $cell.Value2 = $first + $rest
    }
}
